$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = 15.46388812271717
$ws.Range("D2").Value = 4.4530821209041
$ws.Range("E2").Value = 16.80072893514019
$ws.Range("F2").Value = 31.66688530865378
$ws.Range("G2").Value = 47.97825748599407
$ws.Range("H2").Value = 14.99599398427061
$ws.Range("I2").Value = 42.00370031271751
$ws.Range("B3").Value = 14.62529716494162
$ws.Range("D3").Value = 4.256110385767948
$ws.Range("E3").Value = 15.81414347791596
$ws.Range("F3").Value = 30.1364997654836
$ws.Range("G3").Value = 45.07585997036803
$ws.Range("H3").Value = 14.56739984798534
$ws.Range("I3").Value = 39.58439637612544
$ws.Range("B4").Value = 14.08227189835594
$ws.Range("D4").Value = 4.129921971931488
$ws.Range("E4").Value = 15.18345797957187
$ws.Range("F4").Value = 29.16934848132706
$ws.Range("G4").Value = 43.21230408299177
$ws.Range("H4").Value = 14.30402308097859
$ws.Range("I4").Value = 38.02187358361297
$ws.Range("B5").Value = 13.85402509191464
$ws.Range("D5").Value = 4.077210539987413
$ws.Range("E5").Value = 14.92044224378858
$ws.Range("F5").Value = 28.76891537651177
$ws.Range("G5").Value = 42.43304870019238
$ws.Range("H5").Value = 14.19681494693921
$ws.Range("I5").Value = 37.36626242879553
$ws.Range("B6").Value = 13.81570835831313
$ws.Range("D6").Value = 4.06838101841469
$ws.Range("E6").Value = 14.87641457860602
$ws.Range("F6").Value = 28.70206186205763
$ws.Range("G6").Value = 42.30247733761876
$ws.Range("H6").Value = 14.17902561583299
$ws.Range("I6").Value = 37.25627603765879
$ws.Range("B7").Value = 14.07922169300461
$ws.Range("D7").Value = 4.129216260298202
$ws.Range("E7").Value = 15.17993479371136
$ws.Range("F7").Value = 29.16397279584906
$ws.Range("G7").Value = 43.20187414878671
$ws.Range("H7").Value = 14.30257651222966
$ws.Range("I7").Value = 38.0131074349753
$ws.Range("B8").Value = 15.18063999774325
$ws.Range("D8").Value = 4.386272154184288
$ws.Range("E8").Value = 16.46586561567118
$ws.Range("F8").Value = 31.14524790463039
$ws.Range("G8").Value = 46.99476461840648
$ws.Range("H8").Value = 14.84836263646827
$ws.Range("I8").Value = 41.18578900812582
$ws.Range("B9").Value = 17.11354650569209
$ws.Range("D9").Value = 4.847626781533176
$ws.Range("E9").Value = 18.96963095844109
$ws.Range("F9").Value = 34.79019700568814
$ws.Range("G9").Value = 53.76441982954865
$ws.Range("H9").Value = 15.90964749827232
$ws.Range("I9").Value = 46.77879708202152
$ws.Range("B10").Value = 18.39153772040988
$ws.Range("D10").Value = 5.159375015272185
$ws.Range("E10").Value = 20.7005264712871
$ws.Range("F10").Value = 37.29712272872882
$ws.Range("G10").Value = 58.30908890538318
$ws.Range("H10").Value = 16.67484774318615
$ws.Range("I10").Value = 50.48748964831368
$ws.Range("B11").Value = 18.94157951104749
$ws.Range("D11").Value = 5.295119806458914
$ws.Range("E11").Value = 21.44630217105968
$ws.Range("F11").Value = 38.39657939506925
$ws.Range("G11").Value = 60.28053908437305
$ws.Range("H11").Value = 17.01811053947613
$ws.Range("I11").Value = 52.08534703260526
$ws.Range("B12").Value = 19.14533481162755
$ws.Range("D12").Value = 5.345639535607845
$ws.Range("E12").Value = 21.72277565408117
$ws.Range("F12").Value = 38.8067745188044
$ws.Range("G12").Value = 61.01311911472015
$ws.Range("H12").Value = 17.14727530363974
$ws.Range("I12").Value = 52.67744494452594
$ws.Range("B13").Value = 19.10165446781605
$ws.Range("D13").Value = 5.334798663661699
$ws.Range("E13").Value = 21.66349534754383
$ws.Range("F13").Value = 38.71870855538847
$ws.Range("G13").Value = 60.85596817346871
$ws.Range("H13").Value = 17.11949578527255
$ws.Range("I13").Value = 52.55050504930288
$ws.Range("B14").Value = 18.95843357864374
$ws.Range("D14").Value = 5.299293861980987
$ws.Range("E14").Value = 21.46916648747204
$ws.Range("F14").Value = 38.43045085574376
$ws.Range("G14").Value = 60.34108956117144
$ws.Range("H14").Value = 17.02875414856346
$ws.Range("I14").Value = 52.13432002565614
$ws.Range("B15").Value = 18.87011552320613
$ws.Range("D15").Value = 5.277430761849357
$ws.Range("E15").Value = 21.34936291773622
$ws.Range("F15").Value = 38.25307759903486
$ws.Range("G15").Value = 60.02388868579952
$ws.Range("H15").Value = 16.97306174473784
$ws.Range("I15").Value = 51.87770123785438
$ws.Range("B16").Value = 18.35495817582653
$ws.Range("D16").Value = 5.150380618130106
$ws.Range("E16").Value = 20.65095500525945
$ws.Range("F16").Value = 37.2244235607972
$ws.Range("G16").Value = 58.17830377378787
$ws.Range("H16").Value = 16.65230639989389
$ws.Range("I16").Value = 50.38125698588418
$ws.Range("B17").Value = 18.03087938850337
$ws.Range("D17").Value = 5.070875926456215
$ws.Range("E17").Value = 20.21187991088275
$ws.Range("F17").Value = 36.58269053228319
$ws.Range("G17").Value = 57.02138715638986
$ws.Range("H17").Value = 16.45420217812308
$ws.Range("I17").Value = 49.44026495104941
$ws.Range("B18").Value = 17.84153069387962
$ws.Range("D18").Value = 5.024575959885126
$ws.Range("E18").Value = 19.95541536592146
$ws.Range("F18").Value = 36.20973971098365
$ws.Range("G18").Value = 56.34693520985559
$ws.Range("H18").Value = 16.33981001259289
$ws.Range("I18").Value = 48.89063678530209
$ws.Range("B19").Value = 17.77691535278062
$ws.Range("D19").Value = 5.008801897977272
$ws.Range("E19").Value = 19.86790610021658
$ws.Range("F19").Value = 36.08281351125076
$ws.Range("G19").Value = 56.11703346922739
$ws.Range("H19").Value = 16.30100592805423
$ws.Range("I19").Value = 48.70310536761567
$ws.Range("B20").Value = 18.065683389242
$ws.Range("D20").Value = 5.079398517153265
$ws.Range("E20").Value = 20.25902557245064
$ws.Range("F20").Value = 36.65140384155789
$ws.Range("G20").Value = 57.14547825128386
$ws.Range("H20").Value = 16.47533801822657
$ws.Range("I20").Value = 49.54130510658042
$ws.Range("B21").Value = 19.00062425561369
$ws.Range("D21").Value = 5.309746544556469
$ws.Range("E21").Value = 21.52640628333927
$ws.Range("F21").Value = 38.51528771957199
$ws.Range("G21").Value = 60.49270200548502
$ws.Range("H21").Value = 17.0554304110521
$ws.Range("I21").Value = 52.25691677039591
$ws.Range("B22").Value = 19.58522931221873
$ws.Range("D22").Value = 5.455137147564349
$ws.Range("E22").Value = 22.32013858859541
$ws.Range("F22").Value = 39.69753703054367
$ws.Range("G22").Value = 62.59886043891643
$ws.Range("H22").Value = 17.42971607441627
$ws.Range("I22").Value = 53.95605091192146
$ws.Range("B23").Value = 19.27564066207572
$ws.Range("D23").Value = 5.378014062561288
$ws.Range("E23").Value = 21.89965602861509
$ws.Range("F23").Value = 39.06990661277872
$ws.Range("G23").Value = 61.4822596397713
$ws.Range("H23").Value = 17.23043429362993
$ws.Range("I23").Value = 53.05615300520144
$ws.Range("B24").Value = 18.0499579405775
$ws.Range("D24").Value = 5.075547295935355
$ws.Range("E24").Value = 20.23772358888641
$ws.Range("F24").Value = 36.62035102466096
$ws.Range("G24").Value = 57.08940569978319
$ws.Range("H24").Value = 16.46578405071612
$ws.Range("I24").Value = 49.49565175058861
$ws.Range("B25").Value = 16.61539234780988
$ws.Range("D25").Value = 4.727527600169029
$ws.Range("E25").Value = 18.29498466018679
$ws.Range("F25").Value = 33.83266512629792
$ws.Range("G25").Value = 52.00717024423898
$ws.Range("H25").Value = 15.62447975441807
$ws.Range("I25").Value = 45.33509062920456
